$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$changes = @{
    21  = 797
    26  = 951
    37  = 916
    38  = 669
    47  = 504
    62  = 860
    66  = 341
    70  = 1051
    72  = 861
    80  = 981
    95  = 966
    107 = 926
    109 = 743
    134 = 633
    138 = 738
    139 = 826
    140 = 912
    141 = 962
    142 = 985
    144 = 1031
}

foreach ($row in $changes.Keys) {
    $ws.Cells.Item($row, 8).Value = $changes[$row]
}
